# Data update from DGS's 2021/09/10 report.
# Adds a new row (79) to the time series with the latest report date and values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 79
$prevRow = $newRow - 1

# Copy formatting (number formats / styles) from the previous row so the new
# row matches the rest of the table (date-like text in col A, 2 decimals elsewhere).
$ws.Range("A$prevRow`:E$prevRow").Copy() | Out-Null
$ws.Range("A$newRow`:E$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Column A holds the report date as literal text (e.g. "2021/09/08"), not a
# real date serial, so force text entry by prefixing the value assignment
# with a leading apostrophe equivalent via NumberFormat "@" + explicit string.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2021/09/10"

$ws.Cells.Item($newRow, 2).Value = 240.7
$ws.Cells.Item($newRow, 3).Value = 247.9
$ws.Cells.Item($newRow, 4).Value = 0.87
$ws.Cells.Item($newRow, 5).Value = 0.87

# Move the active selection to the next empty row, as Excel would after data entry.
$ws.Range("A80").Select() | Out-Null
